# Auto-generated Excel COM-interop script applying the market-data refresh diff
# (scheduled runner update to Sheets/Halicarnassus_Profits.xlsx)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 363.5
$ws.Range("I2").Value = 222.27272
$ws.Range("K2").Value = 222.27272
$ws.Range("M2").Value = -109.27272

$ws.Range("H15").Value = 1190.5
$ws.Range("I15").Value = 1190.5
$ws.Range("K15").Value = 3571.5
$ws.Range("M15").Value = -3402.5

$ws.Range("H20").Value = 6454.1665
$ws.Range("I20").Value = 5639.8
$ws.Range("J20").Value = 10526
$ws.Range("K20").Value = 5639.8
$ws.Range("L20").Value = 10526
$ws.Range("M20").Value = -5409.8
$ws.Range("N20").Value = -10986

$ws.Range("H35").Value = 6454.1665
$ws.Range("I35").Value = 5639.8
$ws.Range("J35").Value = 10526
$ws.Range("K35").Value = 5639.8
$ws.Range("L35").Value = 10526
$ws.Range("M35").Value = -5260.8
$ws.Range("N35").Value = -11284

$ws.Range("H62").Value = 7837.1665
$ws.Range("I62").Value = 3775
$ws.Range("J62").Value = 8649.6
$ws.Range("K62").Value = 3775
$ws.Range("L62").Value = 8649.6
$ws.Range("M62").Value = -3151
$ws.Range("N62").Value = -9897.6

$ws.Range("H65").Value = 7837.1665
$ws.Range("I65").Value = 3775
$ws.Range("J65").Value = 8649.6
$ws.Range("K65").Value = 18875
$ws.Range("L65").Value = 43248
$ws.Range("M65").Value = -15755
$ws.Range("N65").Value = -49488

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4814.4414
$ws.Range("I32").Value = 4814.4414
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4814.4414
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4527.4414
$ws.Range("N32").Value = ""

$ws.Range("H45").Value = 3266.4443
$ws.Range("I45").Value = 2719.7334
$ws.Range("K45").Value = 2719.7334
$ws.Range("M45").Value = -2342.7334

$ws.Range("H55").Value = 16000

$ws.Range("H110").Value = 3741.5
$ws.Range("I110").Value = 1733
$ws.Range("J110").Value = 5750
$ws.Range("K110").Value = 1733
$ws.Range("L110").Value = 5750
$ws.Range("M110").Value = 312
$ws.Range("N110").Value = -9840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 1183.3334
$ws.Range("I16").Value = 1183.3334
$ws.Range("K16").Value = 1183.3334
$ws.Range("M16").Value = -1013.3334

$ws.Range("H22").Value = 420.57144
$ws.Range("I22").Value = 311.54544
$ws.Range("J22").Value = 820.3333
$ws.Range("K22").Value = 311.54544
$ws.Range("L22").Value = 820.3333
$ws.Range("M22").Value = -138.54544
$ws.Range("N22").Value = -1166.3333

$ws.Range("H62").Value = 100000
$ws.Range("J62").Value = 100000
$ws.Range("L62").Value = 100000
$ws.Range("N62").Value = -101372

$ws.Range("H65").Value = 100000
$ws.Range("J65").Value = 100000
$ws.Range("L65").Value = 300000
$ws.Range("N65").Value = -306864

$ws.Range("H105").Value = 2496.3333
$ws.Range("I105").Value = 2496.3333
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2496.3333
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -749.3332999999998
$ws.Range("N105").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2000
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1713
$ws.Range("N16").Value = ""

$ws.Range("H31").Value = 5719.9707
$ws.Range("I31").Value = 2458.4736
$ws.Range("J31").Value = 9851.200000000001
$ws.Range("K31").Value = 2458.4736
$ws.Range("L31").Value = 9851.200000000001
$ws.Range("M31").Value = -2163.4736
$ws.Range("N31").Value = -10441.2

$ws.Range("H34").Value = 5719.9707
$ws.Range("I34").Value = 2458.4736
$ws.Range("J34").Value = 9851.200000000001
$ws.Range("K34").Value = 2458.4736
$ws.Range("L34").Value = 9851.200000000001
$ws.Range("M34").Value = -2256.4736
$ws.Range("N34").Value = -10255.2

$ws.Range("H53").Value = 66500
$ws.Range("J53").Value = 66500
$ws.Range("L53").Value = 66500
$ws.Range("N53").Value = -67714

$ws.Range("H55").Value = 25000
$ws.Range("I55").Value = 25000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 25000
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -24685
$ws.Range("N55").Value = ""

$ws.Range("H94").Value = 5034.25
$ws.Range("I94").Value = 2065.25
$ws.Range("K94").Value = 2065.25
$ws.Range("M94").Value = -1614.25

$ws.Range("H99").Value = 3904.889
$ws.Range("I99").Value = 3664.0588
$ws.Range("J99").Value = 7999
$ws.Range("K99").Value = 3664.0588
$ws.Range("L99").Value = 7999
$ws.Range("M99").Value = -2166.0588
$ws.Range("N99").Value = -10995

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = ""

$ws.Range("H126").Value = 3904.889
$ws.Range("I126").Value = 3664.0588
$ws.Range("J126").Value = 7999
$ws.Range("K126").Value = 10992.1764
$ws.Range("L126").Value = 23997
$ws.Range("M126").Value = -8522.1764
$ws.Range("N126").Value = -28937

$ws.Range("H141").Value = 76775.336
$ws.Range("J141").Value = 90163
$ws.Range("L141").Value = 90163
$ws.Range("N141").Value = -100523

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 4861.5713
$ws.Range("J114").Value = 4838.5
$ws.Range("L114").Value = 14515.5
$ws.Range("N114").Value = -21023.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 27326.666
$ws.Range("J15").Value = 19990
$ws.Range("L15").Value = 19990
$ws.Range("N15").Value = -20566

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""

$ws.Range("H80").Value = 3135
$ws.Range("J80").Value = 3424
$ws.Range("L80").Value = 3424
$ws.Range("N80").Value = -5420

$ws.Range("H81").Value = 27326.666
$ws.Range("J81").Value = 19990
$ws.Range("L81").Value = 19990
$ws.Range("N81").Value = -21986

$ws.Range("H83").Value = 3135
$ws.Range("J83").Value = 3424
$ws.Range("L83").Value = 17120
$ws.Range("N83").Value = -27104

$ws.Range("H84").Value = 27326.666
$ws.Range("J84").Value = 19990
$ws.Range("L84").Value = 59970
$ws.Range("N84").Value = -69954

$ws.Range("H92").Value = 7107.143
$ws.Range("J92").Value = 7708.3335
$ws.Range("L92").Value = 7708.3335
$ws.Range("N92").Value = -11452.3335

$ws.Range("H97").Value = 976.4167
$ws.Range("I97").Value = 922.2
$ws.Range("K97").Value = 922.2
$ws.Range("M97").Value = -426.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1268.9565
$ws.Range("I55").Value = 1676.1428
$ws.Range("J55").Value = 1090.8125
$ws.Range("K55").Value = 1676.1428
$ws.Range("L55").Value = 1090.8125
$ws.Range("M55").Value = -1503.1428
$ws.Range("N55").Value = -1436.8125

$ws.Range("H61").Value = 3125.3704
$ws.Range("I61").Value = 2277.6086
$ws.Range("K61").Value = 2277.6086
$ws.Range("M61").Value = -2075.6086

$ws.Range("H82").Value = 5829
$ws.Range("I82").Value = 4666.6665
$ws.Range("J82").Value = 6991.3335
$ws.Range("K82").Value = 4666.6665
$ws.Range("L82").Value = 6991.3335
$ws.Range("M82").Value = -4305.6665
$ws.Range("N82").Value = -7713.3335

$ws.Range("H85").Value = 5829
$ws.Range("I85").Value = 4666.6665
$ws.Range("J85").Value = 6991.3335
$ws.Range("K85").Value = 4666.6665
$ws.Range("L85").Value = 6991.3335
$ws.Range("M85").Value = -3418.6665
$ws.Range("N85").Value = -9487.333500000001

$ws.Range("H113").Value = 3125.3704
$ws.Range("I113").Value = 2277.6086
$ws.Range("K113").Value = 2277.6086
$ws.Range("M113").Value = -107.6086

$ws.Range("H132").Value = 3548.25
$ws.Range("I132").Value = 3697.5
$ws.Range("J132").Value = 3399
$ws.Range("K132").Value = 11092.5
$ws.Range("L132").Value = 10197
$ws.Range("M132").Value = -8562.5
$ws.Range("N132").Value = -15257

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 7083333.5
$ws.Range("J94").Value = 7083333.5
$ws.Range("L94").Value = 7083333.5
$ws.Range("N94").Value = -7085135.5

$ws.Range("H98").Value = 42289.8
$ws.Range("J98").Value = 42289.8
$ws.Range("L98").Value = 42289.8
$ws.Range("N98").Value = -48279.8

$ws.Range("H100").Value = 471.5
$ws.Range("I100").Value = 314.4
$ws.Range("K100").Value = 628.8
$ws.Range("M100").Value = -87.79999999999995

$ws.Range("H113").Value = 432.5
$ws.Range("I113").Value = 299.30768
$ws.Range("J113").Value = 778.8
$ws.Range("K113").Value = 897.92304
$ws.Range("L113").Value = 2336.4
$ws.Range("M113").Value = 1272.07696
$ws.Range("N113").Value = -6676.4

$ws.Range("H126").Value = 5210.364
$ws.Range("I126").Value = 2885.6667
$ws.Range("K126").Value = 8657.000100000001
$ws.Range("M126").Value = -6187.000100000001

$ws.Range("H132").Value = 3386.1072
$ws.Range("I132").Value = 3110.1428
$ws.Range("J132").Value = 4214
$ws.Range("K132").Value = 9330.428400000001
$ws.Range("L132").Value = 12642
$ws.Range("M132").Value = -6800.428400000001
$ws.Range("N132").Value = -17702

$ws.Range("H136").Value = 4369.48
$ws.Range("J136").Value = 5283.8667
$ws.Range("L136").Value = 15851.6001
$ws.Range("N136").Value = -20951.6001
